$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header label fixups: replace mangled euro-sign text with "(EUR)" wording
$ws.Range("C1").Value = "Cost per unit (EUR)"
$ws.Range("D1").Value = "Total cost (EUR)"

# Row 3 (Pin strip no ejector 16 pins): reference number swap
$ws.Range("F3").Value = 1580994

# Row 5 (Pin strip no ejector 6 pins): reference number swap + new web link
$ws.Range("F5").Value = 10120550
$ws.Range("G5").Value = "https://www.conrad.com/p/bkl-electronic-10120550-pin-strip-no-ejector-contact-spacing-254-mm-total-number-of-pins-6-no-of-rows-2-1-pcs-741435"
